$d = $word.ActiveDocument

# --- Change 1: "unlikely to recognize" -> "unlikely to visually recognize" ---
$d.Content.Find.Execute(
    "unlikely to recognize and are listed below",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "unlikely to visually recognize and are listed below", 2) | Out-Null

# --- Change 2: "required of each" -> "required for each" ---
$d.Content.Find.Execute(
    "required of each of the above flowers",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "required for each of the above flowers", 2) | Out-Null

# --- Change 3: big restructure of "Importantly..." / "Therefore..." paragraphs ---

# 3a: capitalize "python's" -> "Python's"
$d.Content.Find.Execute(
    "Importantly, python" + [char]8217 + "s Anaconda",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Importantly, Python" + [char]8217 + "s Anaconda", 2) | Out-Null

# 3b: merge the "Importantly..." paragraph with the "Therefore..." paragraph
#     by deleting the paragraph mark between them
$paraImportantly = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Importantly*") {
        $paraImportantly = $p
    }
}
$markStart = $paraImportantly.Range.End - 1
$markEnd = $paraImportantly.Range.End
$d.Range($markStart, $markEnd).Delete()

# 3c: insert the new sentences between "...the program." and "Therefore, ..."
$findRange = $d.Content
$findRange.Find.Execute(
    "Therefore, machine learning algorithm implementation",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint = $d.Range($findRange.Start, $findRange.Start)
$newText = " We will use the Mobile Net V2 model for our image recognition. This algorithm has already been pretrained to recognize common flowers, specifically roses, tulips, daisies, and dandelions. This model is available in the Tensorflow API and requires a single square red-blue-green image as an input. Additionally, the model only accepts JPG/JPEG images. "
$insertPoint.InsertBefore($newText)

# 3d: split back into a new paragraph right before "Therefore"
$findRange2 = $d.Content
$findRange2.Find.Execute(
    "Therefore, machine learning algorithm implementation",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$breakPoint = $d.Range($findRange2.Start, $findRange2.Start)
$breakPoint.InsertParagraphBefore()

# --- Change 4: delete the paragraph "This is what your initial approach..." ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*initial approach to the solution*") {
        $p.Range.Delete()
        break
    }
}

# --- Change 5: "receives input from the user" -> "receives an image input from the user" ---
$d.Content.Find.Execute(
    "program receives input from the user",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "program receives an image input from the user", 2) | Out-Null

Write-Host "All edits applied"
